$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A25").Value = "Record"
$ws.Range("B25").Value = "Balanço Geral"
$ws.Range("C25").Value = "Infraestrutura"
$ws.Range("D25").Value = "2025-04-01T13:14"
$ws.Range("E25").Value = "Negativo"
$ws.Range("F25").Value = "Buracos e falta de iluminação em Estrada de Balança Rangel revolta moradores.  *sem nota da prefeitura*"
